$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("neg_reaction3")

$values = @(2,3,4,5,7,8,9,10,11,12,13,14,15,16,17,19,20,21,22,23,24,26,27,28,30,31,32,33,35,36,37,38,39,40,41,42,43,44,46,47,48,49,51,52,54,56,57,58,59,60,62,63,65,66,67,68,70,71,72,73,74,75)

$col = 1
foreach ($val in $values) {
    $ws.Cells.Item(1, $col).Value = $val
    $col = $col + 1
}
